$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 669

$ws.Range("C3").Value = 947
$ws.Range("D3").Value = 1052
$ws.Range("E3").Value = 1052

$ws.Range("C4").Value = 20.68
$ws.Range("E4").Value = 27.41

$ws.Range("B6").Value = 2.78
$ws.Range("G6").Value = 3.82

$ws.Range("B7").Value = 2.78
$ws.Range("G7").Value = 3.84

$ws.Range("G8").Value = 6.15

$ws.Range("G9").Value = 6.26

$ws.Range("C10").Value = 4399.98
$ws.Range("D10").Value = 7977.53
$ws.Range("E10").Value = 7899.25
$ws.Range("F10").Value = 12080.02
$ws.Range("G10").Value = 12999.98

$ws.Range("C11").Value = 4397.81
$ws.Range("D11").Value = 7977.51
$ws.Range("E11").Value = 7899.52
$ws.Range("F11").Value = 12082.21
$ws.Range("G11").Value = 12999.99

$ws.Range("F18").Value = 0.34

$ws.Range("B22").Value = 7842.39
$ws.Range("C22").Value = 9308.46
$ws.Range("D22").Value = 11633.11
$ws.Range("E22").Value = 11581.2
$ws.Range("F22").Value = 14310.84
$ws.Range("G22").Value = 16645.26

$ws.Range("B23").Value = 9800.66
$ws.Range("C23").Value = 12838.64
$ws.Range("D23").Value = 16568.52
$ws.Range("E23").Value = 16462.02
$ws.Range("F23").Value = 20920.92
$ws.Range("G23").Value = 27006.56

$ws.Range("B24").Value = 0.9
$ws.Range("F24").Value = 4.26
$ws.Range("G24").Value = 6.66

$wb.Save()
